# -----------------------------------------------------------------------
# Applies the "add feature importance, conf. matrix top 6" edit:
#  1. Updates the "Top 6" sheet's bottom rows (new RFC C Si N / RFC C Si
#     results pushed in, reordering the bottom of the sorted table).
#  2. Renames "feature importance" -> "feature importance RFC C Si N" and
#     appends a second "Run RFC C Si N 10 times" accuracy table beneath
#     the existing feature-importance table.
#  3. Inserts a brand new "feature importance RFC C Si" sheet (right
#     after the renamed sheet, before "# records") with its own
#     feature-importance table, "Run RFC C Si 10 times" accuracy table,
#     and a small confusion-matrix style accuracy table ("CN").
#  4. Adds a bar chart of the new sheet's feature-importance table.
#  5. Updates the hidden _xlchart defined names to follow the new sheet.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Top 6" sheet: re-rank bottom three rows of the sorted leaderboard
# ---------------------------------------------------------------------
$top6 = $wb.Worksheets.Item("Top 6")

$top6.Range("A41").Value = "XGBoost"
$top6.Range("B41").Value = "XGBoost: C Si"
$top6.Range("C41").Value = 96.8

$top6.Range("A42").Value = "Random Forest Classifier"
$top6.Range("B42").Value = "Random Forest Classifier: C Si N"
$top6.Range("C42").Value = 96.9

$top6.Range("A43").Value = "Random Forest Classifier"
$top6.Range("B43").Value = "Random Forest Classifier: C Si"
$top6.Range("C43").Value = 97.3

# ---------------------------------------------------------------------
# 2. Rename "feature importance" sheet and extend it with a second
#    accuracy table (10 runs of RFC on C Si N).
# ---------------------------------------------------------------------
$fiN = $wb.Worksheets.Item("feature importance")
$fiN.Name = "feature importance RFC C Si N"

$fiN.Columns.Item(1).ColumnWidth = 18.83203125
$fiN.Columns.Item(2).ColumnWidth = 10.6640625

$fiN.Range("A16").Value = "Run RFC C Si N 10 times"
$fiN.Range("B16").Value = "Accuracy"
$fiN.Range("B16").HorizontalAlignment = -4152

$runsN = @(96.9, 96.6, 96.6, 96.9, 96.9, 96.9, 96.9, 96.9, 96.9, 97.2)
for ($i = 0; $i -lt $runsN.Length; $i++) {
    $r = 17 + $i
    $fiN.Range("B$r").Value = $runsN[$i]
}

$fiN.Range("A27").Value = "Average"
$fiN.Range("B27").Formula = "=AVERAGE(B17:B26)"
$fiN.Range("C27").Formula = "=MEDIAN(B17:B26)"
$fiN.Range("D27").Value = "Median"
$fiN.Range("B27").Borders.Item(8).LineStyle = 1
$fiN.Range("B27").Borders.Item(8).Weight = 2

# ---------------------------------------------------------------------
# 3. Insert the new "feature importance RFC C Si" sheet right after it.
# ---------------------------------------------------------------------
$fiSi = $wb.Worksheets.Add($null, $fiN)
$fiSi.Name = "feature importance RFC C Si"

$fiSi.Columns.Item(1).ColumnWidth = 17
$fiSi.Columns.Item(2).ColumnWidth = 14.1640625

$fiSi.Range("A1").Value = "Random Forest Classifer C Si feature importances"

$fiSi.Range("A3").Value = "elemental isotopes"
$fiSi.Range("B3").Value = "contribution (%)"
$fiSi.Range("A3:B3").Font.Underline = 1

$fiSi.Range("A4").Value = "silicon 30 28"
$fiSi.Range("B4").Value = 19.7
$fiSi.Range("A5").Value = "silicon 29 28"
$fiSi.Range("B5").Value = 21.5
$fiSi.Range("A6").Value = "carbon 12 13"
$fiSi.Range("B6").Value = 58.8
$fiSi.Range("B4:B6").NumberFormat = "0.0"

$fiSi.Range("B7").Formula = "=SUM(B4:B6)"
$fiSi.Range("B7").NumberFormat = "0.0"
$fiSi.Range("B7").Borders.Item(8).LineStyle = 1
$fiSi.Range("B7").Borders.Item(8).Weight = 2

$fiSi.Sort.SortFields.Clear()
$fiSi.Sort.SortFields.Add($fiSi.Range("B4:B6"), 0, 1) | Out-Null
$fiSi.Sort.SetRange($fiSi.Range("A4:B6"))
$fiSi.Sort.Header = -4142
$fiSi.Sort.Apply()

$fiSi.Range("A15").Value = "Run RFC C Si 10 times"
$fiSi.Range("B15").Value = "Accuracy"
$fiSi.Range("B15").HorizontalAlignment = -4152

$runsSi = @(97.4, 97.4, 97.3, 97.3, 97.3, 97.3, 97.2, 97.3, 97.3, 97.3)
for ($i = 0; $i -lt $runsSi.Length; $i++) {
    $r = 16 + $i
    $fiSi.Range("B$r").Value = $runsSi[$i]
}

$fiSi.Range("A26").Value = "Average"
$fiSi.Range("B26").Formula = "=AVERAGE(B16:B25)"
$fiSi.Range("C26").Formula = "=MEDIAN(B16:B25)"
$fiSi.Range("D26").Value = "Median"
$fiSi.Range("B26").Borders.Item(8).LineStyle = 1
$fiSi.Range("B26").Borders.Item(8).Weight = 2

$fiSi.Range("H20").Value = "CN"
$cn = @(95.8, 95.8, 95.6, 95.8, 95.6, 95.6, 95.6, 95.8, 95.6, 95.6)
for ($i = 0; $i -lt $cn.Length; $i++) {
    $r = 20 + $i
    $fiSi.Range("I$r").Value = $cn[$i]
}

$fiSi.Range("H30").Value = "Average"
$fiSi.Range("I30").Formula = "=AVERAGE(I20:I29)"
$fiSi.Range("J30").Formula = "=MEDIAN(I20:I29)"
$fiSi.Range("K30").Value = "Median"
$fiSi.Range("I30").Borders.Item(8).LineStyle = 1
$fiSi.Range("I30").Borders.Item(8).Weight = 2

# ---------------------------------------------------------------------
# 4. Chart the new feature-importance table as a clustered bar chart.
# ---------------------------------------------------------------------
$co = $fiSi.ChartObjects().Add(212.611328125, 19, 320.625, 216)
$chart = $co.Chart
$chart.ChartType = 57
$chart.SeriesCollection().NewSeries()
$ser = $chart.SeriesCollection().Item(1)
$ser.Name = "='feature importance RFC C Si'!`$B`$3"
$ser.XValues = $fiSi.Range("A4:A6")
$ser.Values = $fiSi.Range("B4:B6")

# ---------------------------------------------------------------------
# 5. Re-point the hidden _xlchart defined names at the new sheet.
# ---------------------------------------------------------------------
$wb.Names.Item("_xlchart.v2.0").RefersTo = "='feature importance RFC C Si'!`$A`$4:`$A`$6"
$wb.Names.Item("_xlchart.v2.1").RefersTo = "='feature importance RFC C Si'!`$B`$3"
$wb.Names.Item("_xlchart.v2.2").RefersTo = "='feature importance RFC C Si'!`$B`$4:`$B`$6"

# ---------------------------------------------------------------------
# 6. "# records" sheet is no longer the active tab; the new sheet is.
# ---------------------------------------------------------------------
$records = $wb.Worksheets.Item("# records")
$records.Activate()
$fiSi.Activate()
